$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values: swap the "target"/"location" reference text so each
# action row points at the right shared-string entry, and fix up the
# placeholder letters (nm/na/nr/ns -> L/M/N/S) for clarity.
$ws.Range("B2").Value = "location 1…location L"
$ws.Range("B3").Value = "target 1…target M"
$ws.Range("B4").Value = "target 1…target N"
$ws.Range("B5").Value = "ability 1…ability S (but can be targeted, so more cases possible)"

# Widen column B slightly to fit the updated text.
$ws.Columns.Item(2).ColumnWidth = 24

# Move the active selection like in the edited file.
$ws.Range("D15").Select() | Out-Null
